# Apply the "update LAOs results" edit to the Discounted cash flow workbook.
#
# Changes made here (the subset of the original authoring diff that is
# actually expressible through the Excel object model / COM automation):
#   1. Rename the (only) worksheet from "Tables S49-S50" to "Tables S52-S53".
#   2. Update the two table-title cells (A1 and A38) with the new table
#      numbers and new scenario wording (present-day -> baseline / target).
#   3. Move the active selection from E27 to A39 (the sheet was scrolled down
#      and the user had selected A39 when the file was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet.
$ws.Name = "Tables S52-S53"

# 2. Update the two table title / caption cells.
$ws.Range("A1").Value = "Table S52. Discounted cash flow analysis of the production process at the baseline fermentation performance scenario."
$ws.Range("A38").Value = "Table S53. Discounted cash flow analysis of the production process at the target fermentation performance scenario."

# 3. Update the saved selection/active cell to A39 (sheet scrolled to show
# the second table around row 29 onward).
$ws.Range("A39").Select()
